$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that was bumped by one day
# (2023-10-06 -> 2023-10-07, i.e. serial 45205 -> 45206) for every data
# row. Data runs from row 2 to row 432.
$firstRow = 2
$lastRow = 432

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}
